# Add an "Icon" column to the Kanban task table and populate a few cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the existing table ("Tabelle1") by one column; this also expands
# the table's range (A1:D19 -> A1:E19) and the sheet dimension.
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add() | Out-Null

# Name the new column via its header cell (drives the table column name).
$ws.Range("E1").Value = "Icon"

# Populate the new "Icon" column values for the relevant rows.
$ws.Range("E5").Value = "asterix"
$ws.Range("E12").Value = "trophy"
$ws.Range("E13").Value = "comments"
$ws.Range("E14").Value = "wrong"

# Update the active selection to match the edited cell.
$ws.Range("E15").Select() | Out-Null
